$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.861952666666666
$ws.Range("H2").Value = 14.585858
$ws.Range("I2").Value = 0.3995648519435639
$ws.Range("J2").Value = 0.3995648519435638
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 3.617098305893333
$ws.Range("R2").Value = 32.55388475303999
$ws.Range("S2").Value = 0.002500877693636682
$ws.Range("T2").Value = 0.002500877693636682
$ws.Range("G3").Value = 4.861952666666666
$ws.Range("H3").Value = 14.585858
$ws.Range("I3").Value = 0.3995648519435639
$ws.Range("J3").Value = 0.3995648519435638
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 428.5800838674942
$ws.Range("R3").Value = 3857.220754807447
$ws.Range("S3").Value = 0.2963221568888049
$ws.Range("T3").Value = 0.2963221568888048
$ws.Range("G4").Value = 4.861952666666666
$ws.Range("H4").Value = 14.585858
$ws.Range("I4").Value = 0.3995648519435639
$ws.Range("J4").Value = 0.3995648519435638
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 144.7335209116315
$ws.Range("R4").Value = 1302.601688204684
$ws.Range("S4").Value = 0.1000693935743066
$ws.Range("T4").Value = 0.1000693935743066
$ws.Range("G5").Value = 4.861952666666666
$ws.Range("H5").Value = 14.585858
$ws.Range("I5").Value = 0.3995648519435639
$ws.Range("J5").Value = 0.3995648519435638
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 0.9725477364695556
$ws.Range("R5").Value = 8.752929628226001
$ws.Range("S5").Value = 0.0006724237868157303
$ws.Range("T5").Value = 0.0006724237868157303
$ws.Range("I6").Value = 0.04932556406896855
$ws.Range("J6").Value = 0.04932556406896854
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 0.4465242960266667
$ws.Range("R6").Value = 4.01871866424
$ws.Range("S6").Value = 0.0003087288641783585
$ws.Range("T6").Value = 0.0003087288641783585
$ws.Range("I7").Value = 0.04932556406896855
$ws.Range("J7").Value = 0.04932556406896854
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("R7").Value = 476.166981514888
$ws.Range("S7").Value = 0.03658043860358907
$ws.Range("T7").Value = 0.03658043860358906
$ws.Range("I8").Value = 0.04932556406896855
$ws.Range("J8").Value = 0.04932556406896854
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 17.86709347413378
$ws.Range("R8").Value = 160.803841267204
$ws.Range("S8").Value = 0.01235338709118856
$ws.Range("T8").Value = 0.01235338709118856
$ws.Range("I9").Value = 0.04932556406896855
$ws.Range("J9").Value = 0.04932556406896854
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 0.1200592731117778
$ws.Range("R9").Value = 1.080533458006
$ws.Range("S9").Value = 0.00008300951001256359
$ws.Range("T9").Value = 0.00008300951001256359
$ws.Range("G10").Value = 4.206754333333333
$ws.Range("H10").Value = 12.620263
$ws.Range("I10").Value = 0.3457193616641432
$ws.Range("J10").Value = 0.3457193616641432
$ws.Range("M10").Value = 0.74396
$ws.Range("N10").Value = 2.23188
$ws.Range("O10").Value = 0.006259003216804254
$ws.Range("P10").Value = 0.006259003216804255
$ws.Range("Q10").Value = 3.129656953826666
$ws.Range("R10").Value = 28.16691258444
$ws.Range("S10").Value = 0.002163858596767386
$ws.Range("T10").Value = 0.002163858596767386
$ws.Range("G11").Value = 4.206754333333333
$ws.Range("H11").Value = 12.620263
$ws.Range("I11").Value = 0.3457193616641432
$ws.Range("J11").Value = 0.3457193616641432
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("O11").Value = 0.7416121699579786
$ws.Range("P11").Value = 0.7416121699579786
$ws.Range("Q11").Value = 370.824491433403
$ws.Range("R11").Value = 3337.420422900627
$ws.Range("S11").Value = 0.2563896860002325
$ws.Range("T11").Value = 0.2563896860002325
$ws.Range("G12").Value = 4.206754333333333
$ws.Range("H12").Value = 12.620263
$ws.Range("I12").Value = 0.3457193616641432
$ws.Range("J12").Value = 0.3457193616641432
$ws.Range("M12").Value = 29.76859933333333
$ws.Range("N12").Value = 89.305798
$ws.Range("O12").Value = 0.2504459365921425
$ws.Range("P12").Value = 0.2504459365921425
$ws.Range("Q12").Value = 125.2291842427638
$ws.Range("R12").Value = 1127.062658184874
$ws.Range("S12").Value = 0.08658400933001399
$ws.Range("T12").Value = 0.08658400933001399
$ws.Range("G13").Value = 4.206754333333333
$ws.Range("H13").Value = 12.620263
$ws.Range("I13").Value = 0.3457193616641432
$ws.Range("J13").Value = 0.3457193616641432
$ws.Range("M13").Value = 0.2000323333333334
$ws.Range("N13").Value = 0.6000970000000001
$ws.Range("O13").Value = 0.00168289023307462
$ws.Range("P13").Value = 0.00168289023307462
$ws.Range("Q13").Value = 0.8414868850567778
$ws.Range("R13").Value = 7.573381965511
$ws.Range("S13").Value = 0.000581807737129379
$ws.Range("T13").Value = 0.000581807737129379
$ws.Range("G14").Value = 2.499212666666667
$ws.Range("H14").Value = 7.497638
$ws.Range("I14").Value = 0.2053902223233243
$ws.Range("J14").Value = 0.2053902223233243
$ws.Range("M14").Value = 0.74396
$ws.Range("N14").Value = 2.23188
$ws.Range("O14").Value = 0.006259003216804254
$ws.Range("P14").Value = 0.006259003216804255
$ws.Range("Q14").Value = 1.859314255493333
$ws.Range("R14").Value = 16.73382829944
$ws.Range("S14").Value = 0.001285538062221828
$ws.Range("T14").Value = 0.001285538062221828
$ws.Range("G15").Value = 2.499212666666667
$ws.Range("H15").Value = 7.497638
$ws.Range("I15").Value = 0.2053902223233243
$ws.Range("J15").Value = 0.2053902223233243
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("O15").Value = 0.7416121699579786
$ws.Range("P15").Value = 0.7416121699579786
$ws.Range("Q15").Value = 220.3050600690142
$ws.Range("R15").Value = 1982.745540621128
$ws.Range("S15").Value = 0.1523198884653522
$ws.Range("T15").Value = 0.1523198884653522
$ws.Range("G16").Value = 2.499212666666667
$ws.Range("H16").Value = 7.497638
$ws.Range("I16").Value = 0.2053902223233243
$ws.Range("J16").Value = 0.2053902223233243
$ws.Range("M16").Value = 29.76859933333333
$ws.Range("N16").Value = 89.305798
$ws.Range("O16").Value = 0.2504459365921425
$ws.Range("P16").Value = 0.2504459365921425
$ws.Range("Q16").Value = 74.39806052279155
$ws.Range("R16").Value = 669.582544705124
$ws.Range("S16").Value = 0.05143914659663333
$ws.Range("T16").Value = 0.05143914659663333
$ws.Range("G17").Value = 2.499212666666667
$ws.Range("H17").Value = 7.497638
$ws.Range("I17").Value = 0.2053902223233243
$ws.Range("J17").Value = 0.2053902223233243
$ws.Range("M17").Value = 0.2000323333333334
$ws.Range("N17").Value = 0.6000970000000001
$ws.Range("O17").Value = 0.00168289023307462
$ws.Range("P17").Value = 0.00168289023307462
$ws.Range("Q17").Value = 0.4999233412095557
$ws.Range("R17").Value = 4.499310070886001
$ws.Range("S17").Value = 0.0003456491991169473
$ws.Range("T17").Value = 0.0003456491991169474
